# STS IR Bot Performer / Config.xlsx
# Commit: "Change in setting of parameters in process within main and config file."
#
# Net effect of the OOXML diff: on the "Constants" worksheet, the settings
# row that defined "ReviewSheet_WorksheetName" = "Template" (row 30) is
# removed entirely (not just cleared) so every row below it shifts up by
# one, the sheet's used range shrinks by one row, and the two now-unused
# shared strings ("ReviewSheet_WorksheetName" / "Template") drop out of the
# shared string table on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")

# Delete the entire row so everything below shifts up (matches the diff:
# row 30 disappears, rows 31-751 become rows 30-750).
$ws.Rows.Item(30).Delete() | Out-Null

# Reflect the cursor / selection position recorded in the saved file
# (topLeftCell="A7", active cell A31) as closely as the object model
# allows.
$ws.Activate() | Out-Null
$ws.Range("A31").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
